$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.73076433333334
$ws.Range("H2").Value = 53.19229300000001
$ws.Range("I2").Value = 0.004631884691211661
$ws.Range("J2").Value = 0.00463188469121166
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 48.92883024171878
$ws.Range("R2").Value = 440.359472175469
$ws.Range("S2").Value = 0.001192278309368779
$ws.Range("T2").Value = 0.001192278309368779
$ws.Range("G3").Value = 17.73076433333334
$ws.Range("H3").Value = 53.19229300000001
$ws.Range("I3").Value = 0.004631884691211661
$ws.Range("J3").Value = 0.00463188469121166
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 128.4634124785037
$ws.Range("R3").Value = 1156.170712306533
$ws.Range("S3").Value = 0.003130345432109273
$ws.Range("T3").Value = 0.003130345432109272
$ws.Range("G4").Value = 17.73076433333334
$ws.Range("H4").Value = 53.19229300000001
$ws.Range("I4").Value = 0.004631884691211661
$ws.Range("J4").Value = 0.00463188469121166
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 9.220553017282448
$ws.Range("R4").Value = 82.98497715554201
$ws.Range("S4").Value = 0.0002246827751364719
$ws.Range("T4").Value = 0.0002246827751364719
$ws.Range("G5").Value = 17.73076433333334
$ws.Range("H5").Value = 53.19229300000001
$ws.Range("I5").Value = 0.004631884691211661
$ws.Range("J5").Value = 0.00463188469121166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 3.470927143855112
$ws.Range("R5").Value = 31.23834429469601
$ws.Range("S5").Value = 0.00008457817459713718
$ws.Range("T5").Value = 0.00008457817459713717
$ws.Range("I6").Value = 0.9353873458333681
$ws.Range("J6").Value = 0.935387345833368
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 9880.947326121888
$ws.Range("R6").Value = 88928.525935097
$ws.Range("S6").Value = 0.2407750014613208
$ws.Range("T6").Value = 0.2407750014613208
$ws.Range("I7").Value = 0.9353873458333681
$ws.Range("J7").Value = 0.935387345833368
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.6321585489461611
$ws.Range("T7").Value = 0.632158548946161
$ws.Range("I8").Value = 0.9353873458333681
$ws.Range("J8").Value = 0.935387345833368
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 1862.047349821982
$ws.Range("R8").Value = 16758.42614839784
$ws.Range("S8").Value = 0.0453736305413947
$ws.Range("T8").Value = 0.04537363054139468
$ws.Range("I9").Value = 0.9353873458333681
$ws.Range("J9").Value = 0.935387345833368
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 700.9374250683965
$ws.Range("R9").Value = 6308.436825615569
$ws.Range("S9").Value = 0.01708016488449154
$ws.Range("T9").Value = 0.01708016488449154
$ws.Range("G10").Value = 227.2177583333333
$ws.Range("H10").Value = 681.653275
$ws.Range("I10").Value = 0.0593570833501536
$ws.Range("J10").Value = 0.05935708335015359
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 627.0174774414527
$ws.Range("R10").Value = 5643.157296973074
$ws.Range("S10").Value = 0.0152789129487742
$ws.Range("T10").Value = 0.0152789129487742
$ws.Range("G11").Value = 227.2177583333333
$ws.Range("H11").Value = 681.653275
$ws.Range("I11").Value = 0.0593570833501536
$ws.Range("J11").Value = 0.05935708335015359
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 1646.244237556142
$ws.Range("R11").Value = 14816.19813800527
$ws.Range("S11").Value = 0.04011502598089869
$ws.Range("T11").Value = 0.04011502598089869
$ws.Range("G12").Value = 227.2177583333333
$ws.Range("H12").Value = 681.653275
$ws.Range("I12").Value = 0.0593570833501536
$ws.Range("J12").Value = 0.05935708335015359
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 118.1603538230945
$ws.Range("R12").Value = 1063.44318440785
$ws.Range("S12").Value = 0.002879284589364565
$ws.Range("T12").Value = 0.002879284589364564
$ws.Range("G13").Value = 227.2177583333333
$ws.Range("H13").Value = 681.653275
$ws.Range("I13").Value = 0.0593570833501536
$ws.Range("J13").Value = 0.05935708335015359
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 44.47954245731111
$ws.Range("R13").Value = 400.3158821158
$ws.Range("S13").Value = 0.001083859831116142
$ws.Range("T13").Value = 0.001083859831116142
$ws.Range("G14").Value = 2.387458333333333
$ws.Range("H14").Value = 7.162374999999999
$ws.Range("I14").Value = 0.0006236861252666267
$ws.Range("J14").Value = 0.0006236861252666266
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 6.588297114819443
$ws.Range("R14").Value = 59.29467403337499
$ws.Range("S14").Value = 0.0001605410083762549
$ws.Range("T14").Value = 0.0001605410083762549
$ws.Range("G15").Value = 2.387458333333333
$ws.Range("H15").Value = 7.162374999999999
$ws.Range("I15").Value = 0.0006236861252666267
$ws.Range("J15").Value = 0.0006236861252666266
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 17.29767757804166
$ws.Range("R15").Value = 155.679098202375
$ws.Range("S15").Value = 0.0004215029396138958
$ws.Range("T15").Value = 0.0004215029396138957
$ws.Range("G16").Value = 2.387458333333333
$ws.Range("H16").Value = 7.162374999999999
$ws.Range("I16").Value = 0.0006236861252666267
$ws.Range("J16").Value = 0.0006236861252666266
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 1.241553140361111
$ws.Range("R16").Value = 11.17397826325
$ws.Range("S16").Value = 0.00003025367399687184
$ws.Range("T16").Value = 0.00003025367399687183
$ws.Range("G17").Value = 2.387458333333333
$ws.Range("H17").Value = 7.162374999999999
$ws.Range("I17").Value = 0.0006236861252666267
$ws.Range("J17").Value = 0.0006236861252666266
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.4673624767777777
$ws.Range("R17").Value = 4.206262291
$ws.Range("S17").Value = 0.00001138850327960425
$ws.Range("T17").Value = 0.00001138850327960425

Write-Output "Applied 182 cell updates"
